$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new attendance row right below the header (existing records shift
# down by one row) to register the new apprentice "pablo alfonso" who signed
# in without already having a user record.
$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).Value = "pablo"
$ws.Cells.Item(2, 2).Value = "alfonso"
$ws.Cells.Item(2, 3).Value = "Tarjeta de Extranjeria"
$ws.Cells.Item(2, 4).Value = "987321654"
$ws.Cells.Item(2, 5).Value = "ADSO"
$ws.Cells.Item(2, 6).Value = "Tecnologo"
$ws.Cells.Item(2, 7).Value = ""
$ws.Cells.Item(2, 8).Value = "05:35 p. m."
$ws.Cells.Item(2, 9).Value = "A tiempo"

# Re-generated export: refresh the "Curso" / "Fecha" columns for the
# pre-existing attendance rows (now shifted to rows 3-5), and fix the
# document type typo on the last one.
for ($r = 3; $r -le 5; $r++) {
    $ws.Cells.Item($r, 7).Value = ""
    $ws.Cells.Item($r, 8).Value = "05:36 p. m."
}

$ws.Cells.Item(5, 3).Value = "Cedula de Extranjeria"
